$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 90 (kd_vag): 0.43 -> 0
$ws.Range("A90:H90").Value = 0

# Row 138 (t_min_min): -14 -> -47
$ws.Range("A138:H138").Value = -47

# Row 139 (Ta): 73 -> 106
$ws.Range("A139:H139").Value = 106

# Row 143 (t_max_zakr_curve): 59 -> 56
$ws.Range("A143:H143").Value = 56
